$wb = $excel.ActiveWorkbook

# --- Apendix2: add data row for the CMND/paper entry ---
$ws2 = $wb.Worksheets.Item("Apendix2")
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "Nguyễn Hồng Phúc"
$ws2.Range("C2").Value = "HE130001"
$ws2.Range("D2").Value = "FPTUHN2"
$ws2.Range("E2").Value = "paper 9"
$ws2.Range("F2").Value = "ABC"
$ws2.Range("G2").Value = "2 tác giả, 2 địa chỉ FPTU"

# --- Apendix3: add new "CMND" header column ---
$ws3 = $wb.Worksheets.Item("Apendix3")
$ws3.Range("F1").Value = "CMND"
$ws3.Range("F1").Font.Bold = $true

# --- Apendix4: add new "CMND" header column + data row with link text ---
$ws4 = $wb.Worksheets.Item("Apendix4")
$ws4.Range("F1").Value = "CMND"
$ws4.Range("F1").Font.Bold = $true
$ws4.Range("A2").Value = 1
$ws4.Range("B2").Value = "Nguyễn Hồng Phúc"
$ws4.Range("C2").Value = "HE130001"
$ws4.Range("D2").Value = "FPTUHN2"
$ws4.Range("E2").Value = "5.000.000 ₫"
$ws4.Range("F2").Value = "https://www.google.com.vn/?hl=vi"
$ws4.PageSetup.Orientation = 1

# --- Apendix4 selection moves to F1 (even though it's no longer the active tab) ---
[void]$ws4.Activate()
[void]$ws4.Range("F1").Select()

# --- Activate Apendix3 as the selected tab, with F1 selected ---
[void]$ws3.Activate()
[void]$ws3.Range("F1").Select()
